$d = $word.ActiveDocument

# Insert a brand-new empty paragraph right before the current first
# paragraph ("First text paragraph"), then turn it into the document
# title paragraph.
$insertionPoint = $d.Range(0, 0)
$insertionPoint.InsertParagraphBefore()

$titlePara = $d.Paragraphs(1)
$titlePara.Range.Text = " 010 only text paragraph"
$titlePara.Style = "Title"
